$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank row 4 cells (a lead entry that the error
# collector had skipped over): company name, website, industry, and stage.
$ws.Range("A4").Value = "PANAMERICANA"
$ws.Range("B4").Value = "pananmericana"
$ws.Range("C4").Value = "librería"
$ws.Range("H4").Value = "PRE QUALIFICATION"

# Match the ETAPA column styling used elsewhere (same formatting as H2/H3,
# which is the "PRE QUALIFICATION" stage style) by copying formats only.
$ws.Range("H3").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection left over from editing.
$ws.Range("F6").Select()
